$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update of league bases (as per commit): rows 235-237 and 309/310/313 get their
# match-data columns (B, E-AD) rotated between rows, while id/Div/Date (A, C, D)
# stay attached to their original row.

# --- Row 235 (was id 233 / match 6870268) now takes data of old row 236 ---
$ws.Range("B235").Value = 6865915
$ws.Range("E235").Value = "FC Voluntari"
$ws.Range("F235").Value = "Universitatea Cluj"
$ws.Range("G235").Value = 0
$ws.Range("H235").Value = 0
$ws.Range("I235").Value = 0
$ws.Range("J235").Value = 0
$ws.Range("K235").Value = "D"
$ws.Range("L235").Value = 3.5
$ws.Range("M235").Value = 3.25
$ws.Range("N235").Value = 2.05
$ws.Range("O235").Value = 3.4
$ws.Range("P235").Value = 3.1
$ws.Range("Q235").Value = 2.15
$ws.Range("R235").Value = 0.25
$ws.Range("S235").Value = 1.975
$ws.Range("T235").Value = 1.875
$ws.Range("U235").Value = 2.25
$ws.Range("V235").Value = 2.05
$ws.Range("W235").Value = 1.75
$ws.Range("X235").Value = -1
$ws.Range("Y235").Value = 2.1
$ws.Range("Z235").Value = -1
$ws.Range("AA235").Value = 0.4875
$ws.Range("AB235").Value = -0.5
$ws.Range("AC235").Value = -1
$ws.Range("AD235").Value = 0.75

# --- Row 236 (was id 234 / match 6865915) now takes data of old row 237 ---
$ws.Range("B236").Value = 6861095
$ws.Range("E236").Value = "FC Botosani"
$ws.Range("F236").Value = "Farul Constanta"
$ws.Range("G236").Value = 0
$ws.Range("H236").Value = 0
$ws.Range("I236").Value = 0
$ws.Range("J236").Value = 0
$ws.Range("K236").Value = "D"
$ws.Range("L236").Value = 3.75
$ws.Range("M236").Value = 3.4
$ws.Range("N236").Value = 1.909
$ws.Range("O236").Value = 3.1
$ws.Range("P236").Value = 3
$ws.Range("Q236").Value = 2.375
$ws.Range("R236").Value = 0.25
$ws.Range("S236").Value = 1.775
$ws.Range("T236").Value = 2.1
$ws.Range("U236").Value = 2
$ws.Range("V236").Value = 1.8
$ws.Range("W236").Value = 2.05
$ws.Range("X236").Value = -1
$ws.Range("Y236").Value = 2
$ws.Range("Z236").Value = -1
$ws.Range("AA236").Value = 0.3875
$ws.Range("AB236").Value = -0.5
$ws.Range("AC236").Value = -1
$ws.Range("AD236").Value = 1.05

# --- Row 237 (was id 235 / match 6861095) now takes data of old row 235 ---
$ws.Range("B237").Value = 6870268
$ws.Range("E237").Value = "Petrolul Ploiesti"
$ws.Range("F237").Value = "ACS Sepsi"
$ws.Range("G237").Value = 1
$ws.Range("H237").Value = 2
$ws.Range("I237").Value = 0
$ws.Range("J237").Value = 1
$ws.Range("K237").Value = "A"
$ws.Range("L237").Value = 2.8
$ws.Range("M237").Value = 3
$ws.Range("N237").Value = 2.55
$ws.Range("O237").Value = 3
$ws.Range("P237").Value = 3.2
$ws.Range("Q237").Value = 2.3
$ws.Range("R237").Value = 0.25
$ws.Range("S237").Value = 1.85
$ws.Range("T237").Value = 2
$ws.Range("U237").Value = 2.25
$ws.Range("V237").Value = 1.875
$ws.Range("W237").Value = 1.975
$ws.Range("X237").Value = -1
$ws.Range("Y237").Value = -1
$ws.Range("Z237").Value = 1.3
$ws.Range("AA237").Value = -1
$ws.Range("AB237").Value = 1
$ws.Range("AC237").Value = 0.875
$ws.Range("AD237").Value = -1

# --- Row 309 (was id 307 / match 8191523) now takes data of old row 310 ---
$ws.Range("B309").Value = 8191476
$ws.Range("E309").Value = "FC Voluntari"
$ws.Range("F309").Value = "Universitatea Cluj"
$ws.Range("G309").Value = 0
$ws.Range("H309").Value = 1
$ws.Range("I309").Value = 0
$ws.Range("J309").Value = 1
$ws.Range("K309").Value = "A"
$ws.Range("L309").Value = 3.05
$ws.Range("M309").Value = 3.3
$ws.Range("N309").Value = 2.15
$ws.Range("O309").Value = 2.6
$ws.Range("P309").Value = 3.4
$ws.Range("Q309").Value = 2.4
$ws.Range("R309").Value = 0
$ws.Range("S309").Value = 2
$ws.Range("T309").Value = 1.85
$ws.Range("U309").Value = 2.25
$ws.Range("V309").Value = 2
$ws.Range("W309").Value = 1.85
$ws.Range("X309").Value = -1
$ws.Range("Y309").Value = -1
$ws.Range("Z309").Value = 1.4
$ws.Range("AA309").Value = -1
$ws.Range("AB309").Value = 0.8500000000000001
$ws.Range("AC309").Value = -1
$ws.Range("AD309").Value = 0.8500000000000001

# --- Row 310 (was id 308 / match 8191463) now takes data of old row 313 ---
$ws.Range("B310").Value = 8191523
$ws.Range("E310").Value = "Otelul Galati"
$ws.Range("F310").Value = "FC Botosani"
$ws.Range("G310").Value = 2
$ws.Range("H310").Value = 0
$ws.Range("I310").Value = 2
$ws.Range("J310").Value = 0
$ws.Range("K310").Value = "H"
$ws.Range("L310").Value = 1.666
$ws.Range("M310").Value = 3.6
$ws.Range("N310").Value = 4.6
$ws.Range("O310").Value = 2.9
$ws.Range("P310").Value = 3.5
$ws.Range("Q310").Value = 2.2
$ws.Range("R310").Value = 0.25
$ws.Range("S310").Value = 1.85
$ws.Range("T310").Value = 2
$ws.Range("U310").Value = 2.25
$ws.Range("V310").Value = 1.875
$ws.Range("W310").Value = 1.975
$ws.Range("X310").Value = 1.9
$ws.Range("Y310").Value = -1
$ws.Range("Z310").Value = -1
$ws.Range("AA310").Value = 0.8500000000000001
$ws.Range("AB310").Value = -1
$ws.Range("AC310").Value = -0.5
$ws.Range("AD310").Value = 0.4875

# --- Row 313 (was id 311 / match 8191476) now takes data of old row 309 ---
$ws.Range("B313").Value = 8191463
$ws.Range("E313").Value = "Dinamo Bucharest"
$ws.Range("F313").Value = "ACS UTA Batrana Doamna"
$ws.Range("G313").Value = 2
$ws.Range("H313").Value = 0
$ws.Range("I313").Value = 2
$ws.Range("J313").Value = 0
$ws.Range("K313").Value = "H"
$ws.Range("L313").Value = 1.833
$ws.Range("M313").Value = 3.4
$ws.Range("N313").Value = 3.6
$ws.Range("O313").Value = 1.5
$ws.Range("P313").Value = 4.333
$ws.Range("Q313").Value = 5
$ws.Range("R313").Value = -1
$ws.Range("S313").Value = 1.875
$ws.Range("T313").Value = 1.975
$ws.Range("U313").Value = 3
$ws.Range("V313").Value = 2.025
$ws.Range("W313").Value = 1.825
$ws.Range("X313").Value = 0.5
$ws.Range("Y313").Value = -1
$ws.Range("Z313").Value = -1
$ws.Range("AA313").Value = 0.875
$ws.Range("AB313").Value = -1
$ws.Range("AC313").Value = -1
$ws.Range("AD313").Value = 0.825
